# "add Jerry to do" - adds Jerry's TODO notes to each sheet, plus a
# hyperlink + expanded comment on the "My Task" sheet, and updates the
# active sheet/selection state to reflect where editing left off.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # My Task
$ws2 = $wb.Worksheets.Item(2)   # My Opportunity
$ws3 = $wb.Worksheets.Item(3)   # My appointment

# ---------------------------------------------------------------------
# Sheet 1 - "My Task"
# ---------------------------------------------------------------------

# New "Hyperlink" column header
$ws1.Range("G1").Value = "Hyperlink"

# Turn G3 into a hyperlink labeled "detail" (this also creates the
# Hyperlink cell style used by the workbook)
$ws1.Range("G3").Value = "detail"
$ws1.Hyperlinks.Add($ws1.Range("G3"), "https://example.com")

# Expand Jerry's comment in F3 with the extra notes about Angela / the hyperlink
$ws1.Range("F3").Value = "I replied to Ishan that for measured roundtrip, only 3 roundtrips could be regarded as application sequential roundtrips.`nAngela has doubt about request #5. `nSee detail in hyperlink."

# Make row 3 tall enough to show the expanded comment
$ws1.Rows.Item(3).RowHeight = 72

# Narrow column F slightly to make room for the new column G, and size G
$ws1.Columns.Item(6).ColumnWidth = 44
$ws1.Columns.Item(7).ColumnWidth = 27.3

# Jerry's TODO list
$ws1.Range("A8").Value = "TODO(Jerry):"
$ws1.Range("B8").Value = "1. in STEP 1, currently the batch request( technical detail, priority and user status) is sequentially before task list retrieval. Check whether parallem processing of"
$ws1.Range("B9").Value = "both is possible."

# ---------------------------------------------------------------------
# Sheet 3 - "My appointment"
# ---------------------------------------------------------------------

$ws3.Range("A7").Value = "TODO(Jerry):"
$ws3.Range("B7").Value = "1. check whether it is the common procedure to check the edit authorization in a separate roundtrip in the beginning"
$ws3.Range("B8").Value = "2. check Angela's comment on 3 measured roundtrip."

# ---------------------------------------------------------------------
# Sheet 2 - "My Opportunity"
# ---------------------------------------------------------------------

$ws2.Range("A11").Value = "TODO(Jerry):"
$ws2.Range("B11").Value = "1. in STEP 1. check whether it is possible to delay the retrieval of user status and priority until an opportunity is edited."
$ws2.Range("B12").Value = "2. check why there are 8 roundtrips measured when opp is opened."

# ---------------------------------------------------------------------
# Update selections / active sheet to match where editing ended
# ---------------------------------------------------------------------

$ws1.Range("A8").Select()
$ws3.Range("A7").Select()

$ws2.Activate()
$ws2.Range("B12").Select()
